# Debug : TXIQ_CDC & Change : DAC_DISABLE
#
# Updates the Vivado timing-report table on Sheet 1 with the numbers for the
# new critical path (Path 209, txiq_cdc sync_reset_w -> fifo recovery check)
# in place of the old one (Path 191, control_registers field_sdr_reset_reg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name
$ws.Range("B2").Value = "Path 209"

# Slack
$ws.Range("B3").Value = "-5.624ns"

# Source
$ws.Range("B4").Value = "i_system_wrapper/system_i/maia_sdr/inst/txiq_cdc/sync_reset_w/stage1_reg/C   (rising edge-triggered cell FDRE clocked by clk_out1_system_maia_sdr_clk_0  {rise@0.000ns fall@8.000ns period=16.000ns})"

# Requirement
$ws.Range("B8").Value = "0.031ns (rx_clk rise@14464.030ns - clk_out1_system_maia_sdr_clk_0 rise@14464.000ns)"

# Data Path Delay
$ws.Range("B9").Value = "0.838ns (logic 0.456ns (54.437%)  route 0.382ns (45.563%))"

# Clock Path Skew
$ws.Range("B11").Value = "-2.274ns"

# Clock Uncertainty
$ws.Range("B12").Value = "0.175ns"

# Column B grew a little wider to fit the new (longer) Source text.
$ws.Columns.Item(2).ColumnWidth = 205.33
